# Turn the three header lines ("Computer Science" / "10 Marks" /
# "21/02/2021") into labelled, partly-bold lines:
#   Subject: Computer Science
#   Marks: 10
#   Date: 21/02/2021
# where the label ("Subject:", "Marks:", "Date:") is bold and the rest
# of the line is not -- i.e. each paragraph's single run becomes two
# runs.

$d   = $word.ActiveDocument
$hdr = $d.Sections(1).Headers(1)

function Set-LabelledHeaderLine($oldText, $newText, $label) {
    # Step 1: rewrite the whole line's text (still a single run).
    $rng  = $hdr.Range
    $find = $rng.Find
    $find.Execute($oldText, $true, $false, $false, $false, $false, `
                  $true, 1, $false, $newText, 2) | Out-Null

    # Step 2: re-find just the label prefix and bold it -- this splits
    # the run in two: a bold "label:" run and a plain run with the rest.
    $rng2  = $hdr.Range
    $find2 = $rng2.Find
    $find2.Execute($label, $true, $false, $false, $false, $false, `
                   $true, 1, $false, $null, 0) | Out-Null
    $rng2.Bold = $true
}

Set-LabelledHeaderLine "Computer Science" "Subject: Computer Science" "Subject:"
Set-LabelledHeaderLine "10 Marks"         "Marks: 10 "                "Marks:"
Set-LabelledHeaderLine "21/02/2021"       "Date: 21/02/2021"          "Date:"
